$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to remain plain text even for numeric-looking
# values (Excel would otherwise normalize e.g. "247.60" -> "247.6" or "1.00" -> "1").
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "37.078.12"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "2.044.01"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "247.60"
$ws.Range("E5").Value = "  -1.56%  "
$ws.Range("D6").Value = "0.663"
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "56.35"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "0.384"
$ws.Range("E9").Value = "  -0.29%  "
$ws.Range("D10").Value = "0.0781"
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("D11").Value = "0.109"
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("D12").Value = "16.02"
$ws.Range("E12").Value = "  -2.79%  "
$ws.Range("D13").Value = "0.898"
$ws.Range("E13").Value = "  +11.85%  "
$ws.Range("D14").Value = "2.341.63"
$ws.Range("E14").Value = "  -0.30%  "
$ws.Range("D15").Value = "5.71"
$ws.Range("E15").Value = "  +2.67%  "
$ws.Range("D16").Value = "2.044.34"
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("D17").Value = "18.83"
$ws.Range("E17").Value = "  +12.38%  "
$ws.Range("D18").Value = "37.116.57"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("D19").Value = "74.88"
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("D20").Value = "0.0₃0891"
$ws.Range("E20").Value = "  -1.89%  "
$ws.Range("D21").Value = "5.41"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").Value = "237.06"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("E24").Value = "  +4.62%  "
$ws.Range("D25").Value = "171.09"
$ws.Range("E25").Value = "  +1.19%  "
$ws.Range("D26").Value = "9.54"
$ws.Range("E26").Value = "  +3.18%  "
$ws.Range("D27").Value = "2.18"
$ws.Range("E27").Value = "  -8.12%  "
$ws.Range("D28").Value = "20.11"
$ws.Range("E28").Value = "  -0.29%  "
$ws.Range("E29").Value = "  -0.72%  "
$ws.Range("D30").Value = "5.14"
$ws.Range("E30").Value = "  +8.66%  "
$ws.Range("D31").Value = "1.19"
$ws.Range("E31").Value = "  +1.13%  "
$ws.Range("D32").Value = "4.67"
$ws.Range("E32").Value = "  +4.85%  "
$ws.Range("D33").Value = "0.0623"
$ws.Range("E33").Value = "  +0.44%  "
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").Value = "0.0877"
$ws.Range("E34").Value = "  -0.93%  "
$ws.Range("B35").Value = "BinanceUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").Value = "1.88"
$ws.Range("E36").Value = "  +5.86%  "
$ws.Range("D37").Value = "2.27"
$ws.Range("E37").Value = "  +1.38%  "
$ws.Range("E38").Value = "  -1.50%  "
$ws.Range("E39").Value = "  +8.24%  "
$ws.Range("D40").Value = "3.07"
$ws.Range("E40").Value = "  +8.28%  "
$ws.Range("D41").Value = "0.0990"
$ws.Range("E41").Value = "  -8.54%  "
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "99.47"
$ws.Range("E43").Value = "  +2.77%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").Value = "1.15"
$ws.Range("E44").Value = "  +1.50%  "
$ws.Range("D45").Value = "17.29"
$ws.Range("E45").Value = "  -2.52%  "
$ws.Range("D46").Value = "2.41"
$ws.Range("E46").Value = "  -3.19%  "
$ws.Range("D47").Value = "1.286.50"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").Value = "2.86"
$ws.Range("E48").Value = "  -1.49%  "
$ws.Range("D49").Value = "6.81"
$ws.Range("E49").Value = "  +0.68%  "
$ws.Range("D50").Value = "2.225.64"
$ws.Range("E50").Value = "  -0.45%  "
$ws.Range("D51").Value = "44.88"
$ws.Range("E51").Value = "  +1.95%  "

# Restore the default (unstyled) look for column D now that the text is locked in,
# matching the original workbook which had no explicit style on these cells.
$ws.Range("D2:D51").Style = "Normal"

